# Generate Report for Handoff
# Rewrites the Overview / zh-cn / de-de sheets from the 3-file md-caller/md-callee
# fixture ("0a44dd6e...png" / "64304b25...png" / "cc43d0f5...md") to the
# 4-file caller/callee markdown fixture (calleeMd1/calleeMd2/callerMd1/callerMd2).

$wb = $excel.ActiveWorkbook

$urlSrc   = "https://github.com/OpenLocalizationTest/oltest/blob/e4be90e28857f269d56fc0c5c11ee2fbfcbf6c27/e2e/"
$urlZh    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57752d5219a71b568e211fa2c648d97fce2ce894/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$urlDe    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a8081555cea4de2cc981973e3eaa9f9b5d09808/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$hyperlinkColor = 15570276   # BGR encoding of RGB FF6495ED (Cornflower Blue), matches the workbook's custom "HyperLink" cell style

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

function Style-AsDatetime($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = "calleeMd1.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-03-24 23:12:12"

$ws1.Range("A3").Value = "calleeMd2.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-24 23:12:12"

$ws1.Range("A4").Value = "callerMd1.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-03-24 23:12:12"

$ws1.Range("A5").Value = "callerMd2.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-03-24 23:12:12"

$ws1.Hyperlinks.Add($ws1.Range("A2"), ($urlSrc + "calleeMd1.md"), "", "", "calleeMd1.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($urlSrc + "calleeMd2.md"), "", "", "calleeMd2.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), ($urlSrc + "callerMd1.md"), "", "", "callerMd1.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), ($urlSrc + "callerMd2.md"), "", "", "callerMd2.md")

Style-AsHyperlink($ws1.Range("A2"))
Style-AsHyperlink($ws1.Range("A3"))
Style-AsHyperlink($ws1.Range("A4"))
Style-AsHyperlink($ws1.Range("A5"))

Style-AsDatetime($ws1.Range("D5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = "calleeMd1.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-24 23:12:08"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = "Include"
$ws2.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$ws2.Range("A3").Value = "calleeMd2.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-24 23:12:08"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("J3").Value = "Include"
$ws2.Range("K3").Value = "e2e\callerMd1.md"

$ws2.Range("A4").Value = "callerMd1.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-24 23:12:08"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$ws2.Range("J4").Value = "Include"
$ws2.Range("K4").Value = ""

$ws2.Range("A5").Value = "callerMd2.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"
$ws2.Range("E5").Value = "2016-03-24 23:12:08"
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("I5").Value = "e2e\calleeMd1.md"
$ws2.Range("J5").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), ($urlSrc + "calleeMd1.md"), "", "", "calleeMd1.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), ($urlZh + "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"), "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($urlSrc + "calleeMd2.md"), "", "", "calleeMd2.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), ($urlZh + "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"), "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), ($urlSrc + "callerMd1.md"), "", "", "callerMd1.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), ($urlZh + "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"), "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), ($urlSrc + "callerMd2.md"), "", "", "callerMd2.md")
$ws2.Hyperlinks.Add($ws2.Range("D5"), ($urlZh + "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"), "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf")

Style-AsHyperlink($ws2.Range("A2"))
Style-AsHyperlink($ws2.Range("D2"))
Style-AsHyperlink($ws2.Range("A3"))
Style-AsHyperlink($ws2.Range("D3"))
Style-AsHyperlink($ws2.Range("A4"))
Style-AsHyperlink($ws2.Range("D4"))
Style-AsHyperlink($ws2.Range("A5"))
Style-AsHyperlink($ws2.Range("D5"))

Style-AsDatetime($ws2.Range("E5"))
Style-AsDatetime($ws2.Range("H5"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = "calleeMd1.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-24 23:12:12"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("J2").Value = "Include"
$ws3.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$ws3.Range("A3").Value = "calleeMd2.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-24 23:12:12"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("J3").Value = "Include"
$ws3.Range("K3").Value = "e2e\callerMd1.md"

$ws3.Range("A4").Value = "callerMd1.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-24 23:12:12"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$ws3.Range("J4").Value = "Include"

$ws3.Range("A5").Value = "callerMd2.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"
$ws3.Range("E5").Value = "2016-03-24 23:12:12"
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("I5").Value = "e2e\calleeMd1.md"
$ws3.Range("J5").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), ($urlSrc + "calleeMd1.md"), "", "", "calleeMd1.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), ($urlDe + "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"), "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($urlSrc + "calleeMd2.md"), "", "", "calleeMd2.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), ($urlDe + "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"), "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), ($urlSrc + "callerMd1.md"), "", "", "callerMd1.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), ($urlDe + "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"), "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), ($urlSrc + "callerMd2.md"), "", "", "callerMd2.md")
$ws3.Hyperlinks.Add($ws3.Range("D5"), ($urlDe + "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"), "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf")

Style-AsHyperlink($ws3.Range("A2"))
Style-AsHyperlink($ws3.Range("D2"))
Style-AsHyperlink($ws3.Range("A3"))
Style-AsHyperlink($ws3.Range("D3"))
Style-AsHyperlink($ws3.Range("A4"))
Style-AsHyperlink($ws3.Range("D4"))
Style-AsHyperlink($ws3.Range("A5"))
Style-AsHyperlink($ws3.Range("D5"))

Style-AsDatetime($ws3.Range("E5"))
Style-AsDatetime($ws3.Range("H5"))
